$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title / byline / email ---
Replace-Text "Quantum Entanglement: Unveiling the Spooky Action at a Distance" "Unlocking the Mysteries of Our Universe: A Journey Through Physics"
Replace-Text "Kira Bell" "Alex Watson"
Replace-Text "kirabell@springmail" "awatson660@hmail"

# --- Body paragraph 1 ---
Replace-Text "In the realm of quantum mechanics, the perplexing phenomenon of quantum entanglement has captivated the imaginations of scientists and philosophers alike" "Journey with us into the awe-inspiring realm of physics, a science that unravels the mysteries of our universe, from the tiniest atoms to the vast cosmos"

Replace-Text " This enigmatic concept challenges our classical intuition, suggesting that particles can exhibit a profound interconnectedness, regardless of the vast distances separating them" " As we peel back the layers of physical phenomena, we'll explore concepts that govern the everyday world around us and delve into the uncharted territories beyond our immediate perception"

Replace-Text (" The eerie synchronization between entangled particles has been dubbed " + [char]34 + "spooky action at a distance" + [char]34 + " by Albert Einstein, a testament to its unnerving implications. This essay delves into the mind-bending world of quantum entanglement, exploring its profound implications for our understanding of reality, the nature of information, and the very fabric of spacetime") " Through experiments, observations, and theoretical insights, physics has shaped our understanding of the universe, revealing its intricate beauty and the profound implications it holds for our existence"

Replace-Text "The path to comprehending quantum entanglement begins with the concept of superposition, a fundamental principle of quantum mechanics" "As we embark on this intellectual odyssey, we'll traverse the enigmatic realm of quantum mechanics, a realm where particles defy classical intuition, behaving in ways that challenge our notions of reality"

Replace-Text " Unlike their classical counterparts, quantum particles can exist in multiple states simultaneously" " We'll peer into the heart of atoms, uncovering the fundamental building blocks of matter and unlocking the secrets of chemical bonding and reactions"

Replace-Text " This paradoxical behavior, seemingly defying common sense, manifests in the realm of entanglement, where particles become intrinsically linked, sharing their destinies across vast distances. The measurement of one entangled particle instantaneously affects the state of its distant counterpart, irrespective of the intervening space. This nonlocal connection between entangled particles has been empirically confirmed through numerous experiments, reinforcing its enigmatic nature" " We'll embark on a cosmic voyage, exploring the mysteries of the stars, galaxies, and the enigmatic black holes that warp spacetime"

Replace-Text "The implications of quantum entanglement extend far beyond the theoretical realm, holding immense promise for transformative technologies" "Along our journey, we'll marvel at the intricate interconnectedness of physical phenomena, from the dance of subatomic particles to the ebb and flow of cosmic tides"

Replace-Text " Quantum cryptography, harnessing the inherent security of entangled particles, offers unbreakable encryption methods, revolutionizing data security" " We'll witness the profound implications of physics in our everyday lives, from the functioning of electronic devices to the mechanics of motion"

Replace-Text " Quantum computing, exploiting the superposition and entanglement of quantum bits, promises exponential speedups in computation, opening up new frontiers in scientific discovery, drug development, and material design. The potential applications of quantum entanglement are vast and continue to inspire groundbreaking research, pushing the boundaries of human knowledge and technological advancement" " And as we confront the unanswered questions that linger at the frontiers of knowledge, we'll step into the shoes of physicists, engaging in thoughtful discussions and designing experiments that push the boundaries of our understanding"

# --- Summary paragraph ---
Replace-Text "Quantum entanglement, a cornerstone of quantum mechanics, exhibits the uncanny synchronization between particles, regardless of their separation" "Through a fascinating journey into the world of physics, we've explored the enigmatic realm of quantum mechanics, unraveled the secrets of atoms and chemical reactions, and embarked on a cosmic odyssey to comprehend the mysteries of the universe"

Replace-Text " This profound interconnectedness challenges classical intuition and has profound implications for our understanding of reality" " Along the way, we've witnessed the interconnectedness of physical phenomena, its profound implications in our daily lives, and the unanswered questions that beckon us to explore further"

Replace-Text " From fundamental questions about the nature of information to the development of groundbreaking technologies like quantum cryptography and computing, quantum entanglement holds immense promise. As research continues to unravel the intricacies of this enigmatic phenomenon, we stand at the precipice of transformative discoveries, poised to redefine our understanding of the universe and pave the way for unprecedented technological advancements" " Physics isn't just a collection of abstract theories; it's a living, breathing science that continues to shape our understanding of the world we inhabit"

# --- Add the new trailing empty paragraph after the Summary paragraph ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter() | Out-Null

Write-Host "Done text replacements"
